$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Group Size 2")
$ws.Activate()

# Update contribution split values for rows 7-11 (columns E and F).
# The dependent H column formulas ("Ok" / "Not implemented" / "Error")
# recalculate automatically.
$ws.Range("E7").Value = 0.6
$ws.Range("F7").Value = 0.4

$ws.Range("E8").Value = 0.8
$ws.Range("F8").Value = 0.2

$ws.Range("E9").Value = 0.5
$ws.Range("F9").Value = 0.5

$ws.Range("E10").Value = 0.5
$ws.Range("F10").Value = 0.5

$ws.Range("E11").Value = 0.5
$ws.Range("F11").Value = 0.5

# Scroll the view down so row 7 is at the top of the window, then
# select F9 as the active cell (matches the saved view state).
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F9").Select()
